# "update: waking hours vs working hours"
# Append a new table (rows 39-48) listing each city's data-coverage window
# (start / end dates) directly below the existing "Final Cities working with"
# table. Two new cities - "Oakland 1" and "Oakland 2" - are introduced as new
# shared strings. Date values use a new cell style (numFmtId 14, a built-in
# date format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# City names for column A (rows 39-48). Most already exist as shared
# strings; "Oakland 1" / "Oakland 2" are brand new.
$cities = @(
    "Bakersfield",
    "East Palo Alto",
    "Fresno",
    "Oakland 1",
    "Oakland 2",
    "Richmond",
    "San Diego",
    "San Francisco",
    "San Pablo",
    "Stockton"
)

# Coverage start/end dates (Excel serial date numbers) for columns B and C.
$starts = @(43160, 39630, 42186, 39448, 40756, 39904, 42736, 42370, 41456, 41487)
$ends   = @(44166, 43678, 43313, 40483, 41518, 43374, 44166, 44166, 44166, 43313)

$firstRow = 39
$lastRow = 48

# Write the city names first.
for ($i = 0; $i -lt $cities.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 1).Value = $cities[$i]
}

# Seed B39 with the first date value, then apply the built-in short-date
# number format (numFmtId 14) to it once - this creates exactly one new
# cellXfs entry. We then copy that formatting (not the value) across the
# whole B39:C48 block so every date cell shares that single style, instead
# of minting a new style per cell.
$ws.Cells.Item($firstRow, 2).Value = $starts[0]
$ws.Cells.Item($firstRow, 2).NumberFormat = "mm-dd-yy"
$ws.Cells.Item($firstRow, 2).Copy()
$ws.Range("B$firstRow`:C$lastRow").PasteSpecial(-4122)
# Row 42 (Oakland 1) also picked up the date style on D/E, even though they
# stay empty - matches the source edit exactly.
$ws.Range("D42:E42").PasteSpecial(-4122)

# Now fill in the remaining date values (the style is already in place).
for ($i = 0; $i -lt $starts.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 2).Value = $starts[$i]
    $ws.Cells.Item($row, 3).Value = $ends[$i]
}

# Move the view/selection the way the author left it.
$ws.Range("B37").Select()
